$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 232489.28
$ws.Range("I6").Value = 250107.69
$ws.Range("J6").Value = 3450
$ws.Range("K6").Value = 750323.0700000001
$ws.Range("L6").Value = 10350
$ws.Range("M6").Value = -750211.0700000001
$ws.Range("N6").Value = -10574
$ws.Range("H74").Value = 3527.3076
$ws.Range("I74").Value = 3594.125
$ws.Range("J74").Value = 3420.4
$ws.Range("K74").Value = 3594.125
$ws.Range("L74").Value = 3420.4
$ws.Range("M74").Value = -2658.125
$ws.Range("N74").Value = -5292.4
$ws.Range("H77").Value = 3527.3076
$ws.Range("I77").Value = 3594.125
$ws.Range("J77").Value = 3420.4
$ws.Range("K77").Value = 17970.625
$ws.Range("L77").Value = 17102
$ws.Range("M77").Value = -13290.625
$ws.Range("N77").Value = -26462
$ws.Range("H107").Value = 822.7907
$ws.Range("I107").Value = 820.2727
$ws.Range("K107").Value = 820.2727
$ws.Range("M107").Value = 1099.7273
$ws.Range("H112").Value = 3393.1
$ws.Range("I112").Value = 980
$ws.Range("J112").Value = 3464.0735
$ws.Range("K112").Value = 2940
$ws.Range("L112").Value = 10392.2205
$ws.Range("M112").Value = -1832
$ws.Range("N112").Value = -12608.2205
$ws.Range("H121").Value = 1804.1666
$ws.Range("I121").Value = 272.5
$ws.Range("J121").Value = 2110.5
$ws.Range("K121").Value = 817.5
$ws.Range("L121").Value = 6331.5
$ws.Range("M121").Value = 929.5
$ws.Range("N121").Value = -9825.5
$ws.Range("H132").Value = 2653.805
$ws.Range("I132").Value = 2481.2703
$ws.Range("J132").Value = 4249.75
$ws.Range("K132").Value = 7443.8109
$ws.Range("L132").Value = 12749.25
$ws.Range("M132").Value = -4913.8109
$ws.Range("N132").Value = -17809.25
$ws.Range("H137").Value = 1813135.5
$ws.Range("I137").Value = 3087746.2
$ws.Range("J137").Value = 1846.579
$ws.Range("K137").Value = 9263238.600000001
$ws.Range("L137").Value = 5539.737
$ws.Range("M137").Value = -9260688.600000001
$ws.Range("N137").Value = -10639.737
$ws.Range("H138").Value = 3606.3625
$ws.Range("I138").Value = 2824.48
$ws.Range("J138").Value = 3902.5303
$ws.Range("K138").Value = 8473.440000000001
$ws.Range("L138").Value = 11707.5909
$ws.Range("M138").Value = -3333.440000000001
$ws.Range("N138").Value = -21987.5909
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7364013
$ws.Range("I32").Value = 8342781.5
$ws.Range("K32").Value = 8342781.5
$ws.Range("M32").Value = -8342494.5
$ws.Range("H45").Value = 2432.0625
$ws.Range("I45").Value = 1070.6666
$ws.Range("J45").Value = 2746.2307
$ws.Range("K45").Value = 1070.6666
$ws.Range("L45").Value = 2746.2307
$ws.Range("M45").Value = -693.6666
$ws.Range("N45").Value = -3500.2307
$ws.Range("H61").Value = 9807177
$ws.Range("I61").Value = 16668636
$ws.Range("K61").Value = 16668636
$ws.Range("M61").Value = -16668424
$ws.Range("H74").Value = 23812154
$ws.Range("I74").Value = 1569.3
$ws.Range("J74").Value = 45458140
$ws.Range("K74").Value = 1569.3
$ws.Range("L74").Value = 45458140
$ws.Range("M74").Value = -695.3
$ws.Range("N74").Value = -45459888
$ws.Range("H77").Value = 23812154
$ws.Range("I77").Value = 1569.3
$ws.Range("J77").Value = 45458140
$ws.Range("K77").Value = 7846.5
$ws.Range("L77").Value = 227290700
$ws.Range("M77").Value = -3478.5
$ws.Range("N77").Value = -227299436
$ws.Range("H132").Value = 2499.4119
$ws.Range("I132").Value = 1888.1731
$ws.Range("J132").Value = 4485.9375
$ws.Range("K132").Value = 5664.5193
$ws.Range("L132").Value = 13457.8125
$ws.Range("M132").Value = -3134.5193
$ws.Range("N132").Value = -18517.8125
$ws.Range("H136").Value = 9807177
$ws.Range("I136").Value = 16668636
$ws.Range("K136").Value = 50005908
$ws.Range("M136").Value = -50003358
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 49000
$ws.Range("J108").Value = 49000
$ws.Range("L108").Value = 49000
$ws.Range("N108").Value = -56680
$ws.Range("H135").Value = 47695
$ws.Range("J135").Value = 47695
$ws.Range("L135").Value = 47695
$ws.Range("N135").Value = -57835
$ws.Range("H137").Value = 50780
$ws.Range("J137").Value = 50780
$ws.Range("L137").Value = 50780
$ws.Range("N137").Value = -60980
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1247.4286
$ws.Range("I58").Value = 1074.1538
$ws.Range("J58").Value = 3500
$ws.Range("K58").Value = 1074.1538
$ws.Range("L58").Value = 3500
$ws.Range("M58").Value = -871.1538
$ws.Range("N58").Value = -3906
$ws.Range("H135").Value = 60113.332
$ws.Range("J135").Value = 60113.332
$ws.Range("L135").Value = 60113.332
$ws.Range("N135").Value = -70253.33199999999
$ws.Range("H136").Value = 1247.4286
$ws.Range("I136").Value = 1074.1538
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 3222.4614
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -672.4614000000001
$ws.Range("N136").Value = -15600
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 7686.143
$ws.Range("J49").Value = 7686.143
$ws.Range("L49").Value = 23058.429
$ws.Range("N49").Value = -23370.429
$ws.Range("H64").Value = 2938.5925
$ws.Range("I64").Value = 1666.6666
$ws.Range("J64").Value = 3097.5833
$ws.Range("K64").Value = 4999.9998
$ws.Range("L64").Value = 9292.749899999999
$ws.Range("M64").Value = -4729.9998
$ws.Range("N64").Value = -9832.749899999999
$ws.Range("H67").Value = 2938.5925
$ws.Range("I67").Value = 1666.6666
$ws.Range("J67").Value = 3097.5833
$ws.Range("K67").Value = 4999.9998
$ws.Range("L67").Value = 9292.749899999999
$ws.Range("M67").Value = -4063.9998
$ws.Range("N67").Value = -11164.7499
$ws.Range("H100").Value = 8820.362999999999
$ws.Range("J100").Value = 8820.362999999999
$ws.Range("L100").Value = 26461.089
$ws.Range("N100").Value = -28083.089
$ws.Range("H103").Value = 4909.5713
$ws.Range("I103").Value = 482
$ws.Range("J103").Value = 7369.3335
$ws.Range("K103").Value = 1446
$ws.Range("L103").Value = 22108.0005
$ws.Range("M103").Value = -567
$ws.Range("N103").Value = -23866.0005
$ws.Range("H137").Value = 39461.266
$ws.Range("I137").Value = 8243.412
$ws.Range("J137").Value = 80284.62
$ws.Range("K137").Value = 24730.236
$ws.Range("L137").Value = 240853.86
$ws.Range("M137").Value = -19630.236
$ws.Range("N137").Value = -251053.86
$ws.Range("H140").Value = 1502.9773
$ws.Range("I140").Value = 1063.1818
$ws.Range("J140").Value = 2822.3635
$ws.Range("K140").Value = 3189.5454
$ws.Range("L140").Value = 8467.0905
$ws.Range("M140").Value = 1990.4546
$ws.Range("N140").Value = -18827.0905
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H86").Value = 92880.5
$ws.Range("J86").Value = 92880.5
$ws.Range("L86").Value = 92880.5
$ws.Range("N86").Value = -95252.5
$ws.Range("H89").Value = 92880.5
$ws.Range("J89").Value = 92880.5
$ws.Range("L89").Value = 278641.5
$ws.Range("N89").Value = -290497.5
$ws.Range("H97").Value = 2780.353
$ws.Range("I97").Value = 2760.5454
$ws.Range("J97").Value = 2816.6667
$ws.Range("K97").Value = 2760.5454
$ws.Range("L97").Value = 2816.6667
$ws.Range("M97").Value = -2264.5454
$ws.Range("N97").Value = -3808.6667
$ws.Range("H132").Value = 31255410
$ws.Range("I132").Value = 43484264
$ws.Range("J132").Value = 3891.6667
$ws.Range("K132").Value = 130452792
$ws.Range("L132").Value = 11675.0001
$ws.Range("M132").Value = -130450262
$ws.Range("N132").Value = -16735.0001
$ws.Range("H134").Value = 23865.2
$ws.Range("J134").Value = 23865.2
$ws.Range("L134").Value = 71595.60000000001
$ws.Range("N134").Value = -76665.60000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 5602.75
$ws.Range("I93").Value = 8384.5
$ws.Range("J93").Value = 1708.3
$ws.Range("K93").Value = 8384.5
$ws.Range("L93").Value = 1708.3
$ws.Range("M93").Value = -7136.5
$ws.Range("N93").Value = -4204.3
$ws.Range("H94").Value = 76775
$ws.Range("J94").Value = 76775
$ws.Range("L94").Value = 76775
$ws.Range("N94").Value = -78127
$ws.Range("H132").Value = 2922
$ws.Range("I132").Value = 1925.5294
$ws.Range("J132").Value = 4804.222
$ws.Range("K132").Value = 5776.5882
$ws.Range("L132").Value = 14412.666
$ws.Range("M132").Value = -3246.5882
$ws.Range("N132").Value = -19472.666
